$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Track 1")

# Fix typo: submission score for "ma" was mistakenly entered as 55.24, correct value is 66.24.
# Stored as text in the sheet, so force text by assigning a string value.
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "66.24"

# Update the active selection to the edited cell (B4)
$ws.Range("B4").Select()
